$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.349.63'
$ws.Range("E2").Value = '  -6.69%  '

# Row 3
$ws.Range("D3").Value = '2.909.33'
$ws.Range("E3").Value = '  -9.54%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.64'
$ws.Range("E5").Value = '  -10.33%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.07'
$ws.Range("E6").Value = '  -15.22%  '

# Row 7
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("D8").Value = '2.879.77'
$ws.Range("E8").Value = '  -10.21%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.453'
$ws.Range("E9").Value = '  -17.28%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("E10").Value = '  -19.79%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.79'
$ws.Range("E11").Value = '  -13.04%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.423'
$ws.Range("E12").Value = '  -14.97%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '31.16'
$ws.Range("E13").Value = '  -20.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000198'
$ws.Range("E14").Value = '  -18.91%  '

# Row 15
$ws.Range("D15").Value = '3.386.05'
$ws.Range("E15").Value = '  -9.28%  '

# Row 16
$ws.Range("D16").Value = '62.296.11'
$ws.Range("E16").Value = '  -6.67%  '

# Row 17
$ws.Range("E17").Value = '  -5.81%  '

# Row 18
$ws.Range("D18").Value = '2.890.88'
$ws.Range("E18").Value = '  -10.11%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '464.03'
$ws.Range("E19").Value = '  -13.01%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.06'
$ws.Range("E20").Value = '  -15.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("E21").Value = '  -17.00%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.617'
$ws.Range("E22").Value = '  -18.87%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.30'
$ws.Range("E23").Value = '  -20.43%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.79'
$ws.Range("E24").Value = '  -14.17%  '

# Row 25
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.65'
$ws.Range("E25").Value = '  -16.15%  '

# Row 26
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.48'
$ws.Range("E27").Value = '  -22.51%  '

# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.85'
$ws.Range("E28").Value = '  -15.98%  '

# Row 29
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.81'
$ws.Range("E29").Value = '  -17.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '24.07'
$ws.Range("E30").Value = '  -18.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.25%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.03'
$ws.Range("E32").Value = '  -10.96%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.30'
$ws.Range("E33").Value = '  -13.63%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '50.15'
$ws.Range("E34").Value = '  -6.14%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '456.61'
$ws.Range("E35").Value = '  -16.70%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.31'
$ws.Range("E36").Value = '  -18.96%  '

# Row 37
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.59'
$ws.Range("E37").Value = '  -19.49%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0380'
$ws.Range("E38").Value = '  -10.70%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0728'
$ws.Range("E39").Value = '  -15.72%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.110'
$ws.Range("E40").Value = '  -12.62%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.58'
$ws.Range("E41").Value = '  -19.25%  '

# Row 42
$ws.Range("D42").Value = '2.553.19'
$ws.Range("E42").Value = '  -12.47%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.07'
$ws.Range("E44").Value = '  -22.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.213'
$ws.Range("E45").Value = '  -19.62%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '109.11'
$ws.Range("E46").Value = '  -10.54%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0986'
$ws.Range("E47").Value = '  -13.82%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.73'
$ws.Range("E48").Value = '  -19.29%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.01'
$ws.Range("E49").Value = '  -20.92%  '

# Row 50
$ws.Range("D50").Value = '0.0₃0449'
$ws.Range("E50").Value = '  -23.73%  '

# Row 51
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("E51").Value = '  -7.09%  '
